$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SharedString text edits (rich text cells) ---
# A8: "Volume 30   Number  25" -> "Volume 30   Number  26"
$ws.Range("A8").Characters(21, 2).Text = "26"

# C9: "Report Covering the Week  6/19/2023  Through  6/25/2023"
#     -> "Report Covering the Week  6/26/2023  Through  7/2/2023"
$ws.Range("C9").Characters(27, 9).Text = "6/26/2023"
$ws.Range("C9").Characters(47, 9).Text = "7/2/2023"

# --- Fix number formatting for cells that changed type (string "0" -> numeric) ---
$ws.Range("C22").NumberFormat = $ws.Range("F22").NumberFormat
$ws.Range("C27").NumberFormat = $ws.Range("D27").NumberFormat

# --- Numeric cell value updates ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("L15").Value = -9.090909090909
$ws.Range("N15").Value = -85.294117647058
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 9
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 76
$ws.Range("J16").Value = 107
$ws.Range("K16").Value = -28.971962616822
$ws.Range("L16").Value = 15.151515151515
$ws.Range("M16").Value = -60.824742268041
$ws.Range("N16").Value = -92.850423330197
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 187
$ws.Range("K17").Value = -12.834224598930
$ws.Range("L17").Value = 18.115942028985
$ws.Range("M17").Value = -10.928961748633
$ws.Range("N17").Value = -62.785388127853
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -81.818181818181
$ws.Range("J18").Value = 81
$ws.Range("K18").Value = -18.518518518518
$ws.Range("L18").Value = -29.787234042553
$ws.Range("M18").Value = -52.173913043478
$ws.Range("N18").Value = -95.361911454673
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -36.363636363636
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 275
$ws.Range("J19").Value = 260
$ws.Range("K19").Value = 5.769230769230
$ws.Range("L19").Value = 32.850241545893
$ws.Range("M19").Value = -12.420382165605
$ws.Range("N19").Value = -49.355432780847
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 64
$ws.Range("K20").Value = -15.625
$ws.Range("L20").Value = -16.923076923076
$ws.Range("M20").Value = -46
$ws.Range("N20").Value = -95.836545875096
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 4.301075268817
$ws.Range("I21").Value = 647
$ws.Range("J21").Value = 718
$ws.Range("K21").Value = -9.888579387186
$ws.Range("L21").Value = 10.787671232876
$ws.Range("M21").Value = -31.461864406779
$ws.Range("N21").Value = -86.673532440782
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 11
$ws.Range("K22").Value = 175
$ws.Range("L22").Value = 83.333333333333
$ws.Range("M22").Value = -21.428571428571
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 108.333333333333
$ws.Range("F24").Value = 176
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 74.257425742574
$ws.Range("I24").Value = 932
$ws.Range("J24").Value = 642
$ws.Range("K24").Value = 45.171339563862
$ws.Range("L24").Value = 38.690476190476
$ws.Range("M24").Value = 47.468354430379
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -68.75
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = -48.4375
$ws.Range("I25").Value = 293
$ws.Range("J25").Value = 297
$ws.Range("K25").Value = -1.346801346801
$ws.Range("L25").Value = 12.692307692307
$ws.Range("M25").Value = -23.298429319371
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("L26").Value = -18.181818181818
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 33
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = -5.714285714285
$ws.Range("L27").Value = 26.923076923076
$ws.Range("N28").Value = -90.277777777777
$ws.Range("N29").Value = -89.285714285714
